$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$subscript3 = [char]0x2083
$subscript6 = [char]0x2086

# Map of row -> @{ D = newValue (or $null if unchanged); E = newValue (or $null if unchanged) }
$changes = @{
    2  = @{ D = "34.192.38"; E = $null }
    3  = @{ D = "1.784.72"; E = "  +0.90%  " }
    4  = @{ D = $null; E = "  +0.10%  " }
    5  = @{ D = "225.85"; E = "  +1.03%  " }
    6  = @{ D = $null; E = "  +0.60%  " }
    7  = @{ D = $null; E = "  +0.10%  " }
    8  = @{ D = "31.86"; E = "  +0.41%  " }
    9  = @{ D = $null; E = "  +0.83%  " }
    10 = @{ D = $null; E = "  +0.32%  " }
    11 = @{ D = "0.0945"; E = "  +1.20%  " }
    12 = @{ D = $null; E = "  +0.92%  " }
    13 = @{ D = "11.02"; E = "  +0.16%  " }
    14 = @{ D = "1.788.51"; E = "  +1.01%  " }
    15 = @{ D = $null; E = "  +2.70%  " }
    16 = @{ D = "34.118.02"; E = "  +1.29%  " }
    17 = @{ D = $null; E = "  +1.70%  " }
    18 = @{ D = "68.05"; E = "  +2.63%  " }
    19 = @{ D = "246.05"; E = "  +3.56%  " }
    20 = @{ D = "0.0${subscript3}0780"; E = "  +0.79%  " }
    21 = @{ D = "10.98"; E = "  +4.44%  " }
    22 = @{ D = $null; E = "  +0.08%  " }
    23 = @{ D = $null; E = "  +2.56%  " }
    24 = @{ D = $null; E = "  -0.15%  " }
    25 = @{ D = "161.62"; E = "  +1.50%  " }
    26 = @{ D = $null; E = "  +2.91%  " }
    27 = @{ D = "16.32"; E = "  +1.73%  " }
    28 = @{ D = $null; E = "  +2.13%  " }
    29 = @{ D = $null; E = "  +0.17%  " }
    30 = @{ D = $null; E = "  +1.05%  " }
    31 = @{ D = "0.0520"; E = "  +2.07%  " }
    32 = @{ D = $null; E = "  +3.29%  " }
    33 = @{ D = $null; E = "  +4.62%  " }
    34 = @{ D = $null; E = "  +1.01%  " }
    35 = @{ D = "1.445.40"; E = "  +4.89%  " }
    36 = @{ D = "0.656"; E = "  +1.91%  " }
    37 = @{ D = $null; E = "  +9.65%  " }
    38 = @{ D = $null; E = "  +4.18%  " }
    39 = @{ D = $null; E = "  +1.77%  " }
    40 = @{ D = "80.16"; E = "  +3.25%  " }
    41 = @{ D = $null; E = "  +0.62%  " }
    42 = @{ D = "0.925"; E = "  +2.53%  " }
    43 = @{ D = $null; E = "  +0.67%  " }
    44 = @{ D = "13.52"; E = "  -0.09%  " }
    45 = @{ D = $null; E = "  +4.59%  " }
    46 = @{ D = $null; E = "  +1.94%  " }
    47 = @{ D = $null; E = "  -0.64%  " }
    48 = @{ D = $null; E = "  -1.05%  " }
    49 = @{ D = "1.944.84"; E = "  +1.41%  " }
    50 = @{ D = "106.21"; E = "  -0.54%  " }
    51 = @{ D = $null; E = "  +0.12%  " }
}

foreach ($row in $changes.Keys) {
    $entry = $changes[$row]
    if ($null -ne $entry.D) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $entry.D
        $cell.Style = "Normal"
    }
    if ($null -ne $entry.E) {
        $cell = $ws.Range("E$row")
        $cell.NumberFormat = "@"
        $cell.Value = $entry.E
        $cell.Style = "Normal"
    }
}
